$d = $word.ActiveDocument

# --- Change 1: first paragraph gets trailing spaces and a red "(This is a
#     change - Version for branch alternate)" note appended in three runs.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs(1)
$insPoint = $p1.Range.End - 1

$r1 = $d.Range($insPoint, $insPoint)
$r1.InsertAfter([char]0x28 + "This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 192

$pos2 = $insPoint + ($r1.End - $r1.Start)
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter("rsion for branch alternate")
$r2.Font.Color = 192

$pos3 = $pos2 + ($r2.End - $r2.Start)
$r3 = $d.Range($pos3, $pos3)
$r3.InsertAfter([char]0x29)
$r3.Font.Color = 192

# --- Change 2: append a new empty paragraph (shaded, fill F9F9F9) at the
#     very end of the document body.
$endRange = $d.Content
$endRange.Collapse(0)
$shadedParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($shadedParaXml)
